$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 671.05615
$ws.Range("J17").Value = 676.4091
$ws.Range("L17").Value = 2029.2273
$ws.Range("N17").Value = -2365.2273

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1095914.2
$ws.Range("I33").Value = 1162278.8
$ws.Range("K33").Value = 1162278.8
$ws.Range("M33").Value = -1162049.8

# Sheet ALC, row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 33500
$ws.Range("J128").Value = 33500
$ws.Range("L128").Value = 33500
$ws.Range("N128").Value = -43460

# Sheet ALC, row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 52373.332
$ws.Range("J133").Value = 52373.332
$ws.Range("L133").Value = 52373.332
$ws.Range("N133").Value = -62493.332

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2778847.8
$ws.Range("I137").Value = 1251044.8
$ws.Range("K137").Value = 3753134.4
$ws.Range("M137").Value = -3750584.4

# Sheet ARM, row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 4533.1665
$ws.Range("I41").Value = 4533.1665
$ws.Range("K41").Value = 4533.1665
$ws.Range("M41").Value = -4119.1665

# Sheet ARM, row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 865666.7
$ws.Range("J123").Value = 865666.7
$ws.Range("L123").Value = 865666.7
$ws.Range("N123").Value = -875466.7

# Sheet ARM, row 127
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H127").Value = 34405
$ws.Range("J127").Value = 34405
$ws.Range("L127").Value = 34405
$ws.Range("N127").Value = -44325

# Sheet ARM, row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 30170.357
$ws.Range("J135").Value = 30170.357
$ws.Range("L135").Value = 30170.357
$ws.Range("N135").Value = -40310.357

# Sheet BSM, row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 39785
$ws.Range("J132").Value = 39785
$ws.Range("L132").Value = 39785
$ws.Range("N132").Value = -49905

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2079.4138
$ws.Range("I31").Value = 1589.3462
$ws.Range("J31").Value = 6326.6665
$ws.Range("K31").Value = 1589.3462
$ws.Range("L31").Value = 6326.6665
$ws.Range("M31").Value = -1294.3462
$ws.Range("N31").Value = -6916.6665

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2079.4138
$ws.Range("I34").Value = 1589.3462
$ws.Range("J34").Value = 6326.6665
$ws.Range("K34").Value = 1589.3462
$ws.Range("L34").Value = 6326.6665
$ws.Range("M34").Value = -1387.3462
$ws.Range("N34").Value = -6730.6665

# Sheet CRP, row 80
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Sheet CRP, row 83
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Sheet CRP, row 97
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 8995
$ws.Range("J97").Value = 8995
$ws.Range("L97").Value = 8995
$ws.Range("N97").Value = -10977

# Sheet CRP, row 124
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 19994.4
$ws.Range("J124").Value = 19994.4
$ws.Range("L124").Value = 19994.4
$ws.Range("N124").Value = -24904.4

# Sheet CRP, row 130
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 39980
$ws.Range("J130").Value = 39980
$ws.Range("L130").Value = 39980
$ws.Range("N130").Value = -50020

# Sheet CRP, row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 59100
$ws.Range("J135").Value = 59100
$ws.Range("L135").Value = 59100
$ws.Range("N135").Value = -69240

# Sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1054.381
$ws.Range("I34").Value = 605.25
$ws.Range("J34").Value = 1330.7693
$ws.Range("K34").Value = 1815.75
$ws.Range("L34").Value = 3992.3079
$ws.Range("M34").Value = -1731.75
$ws.Range("N34").Value = -4160.3079

# Sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2266.6667
$ws.Range("I39").Value = 600
$ws.Range("J39").Value = 2350
$ws.Range("K39").Value = 1800
$ws.Range("L39").Value = 7050
$ws.Range("M39").Value = -1506
$ws.Range("N39").Value = -7638

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 50614
$ws.Range("I55").Value = 2398.3333
$ws.Range("J55").Value = 69900.266
$ws.Range("K55").Value = 7194.999899999999
$ws.Range("L55").Value = 209700.798
$ws.Range("M55").Value = -7017.999899999999
$ws.Range("N55").Value = -210054.798

# Sheet CUL, row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1285.3572
$ws.Range("I98").Value = 324.25
$ws.Range("J98").Value = 1669.8
$ws.Range("K98").Value = 972.75
$ws.Range("L98").Value = 5009.4
$ws.Range("M98").Value = 525.25
$ws.Range("N98").Value = -8005.4

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 873.5897
$ws.Range("I113").Value = 875
$ws.Range("J113").Value = 872.96295
$ws.Range("K113").Value = 2625
$ws.Range("L113").Value = 2618.88885
$ws.Range("M113").Value = -455
$ws.Range("N113").Value = -6958.888849999999

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 48102.887
$ws.Range("I122").Value = 346.72726
$ws.Range("J122").Value = 54837.73
$ws.Range("K122").Value = 3120.54534
$ws.Range("L122").Value = 493539.57
$ws.Range("M122").Value = -670.5453400000001
$ws.Range("N122").Value = -498439.57

# Sheet CUL, row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2512.121
$ws.Range("I125").Value = 1050
$ws.Range("J125").Value = 2837.037
$ws.Range("K125").Value = 3150
$ws.Range("L125").Value = 8511.110999999999
$ws.Range("M125").Value = 1770
$ws.Range("N125").Value = -18351.111

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 924.2030999999999
$ws.Range("I131").Value = 547.375
$ws.Range("K131").Value = 1642.125
$ws.Range("M131").Value = 3397.875

# Sheet CUL, row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 10000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Sheet CUL, row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2616.6667
$ws.Range("I134").Value = 833.3333
$ws.Range("J134").Value = 4400
$ws.Range("K134").Value = 2499.9999
$ws.Range("L134").Value = 13200
$ws.Range("M134").Value = 2570.0001
$ws.Range("N134").Value = -23340

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 35425.066
$ws.Range("I140").Value = 60106.41
$ws.Range("J140").Value = 5454.857
$ws.Range("K140").Value = 180319.23
$ws.Range("L140").Value = 16364.571
$ws.Range("M140").Value = -175139.23
$ws.Range("N140").Value = -26724.571

# Sheet CUL, row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 8874.6
$ws.Range("I141").Value = 3311.9
$ws.Range("J141").Value = 20000
$ws.Range("K141").Value = 9935.700000000001
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -4755.700000000001
$ws.Range("N141").Value = -70360

# Sheet GSM, row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 34293
$ws.Range("J123").Value = 34293
$ws.Range("L123").Value = 34293
$ws.Range("N123").Value = -39193

# Sheet GSM, row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 29986.666
$ws.Range("J130").Value = 29986.666
$ws.Range("L130").Value = 29986.666
$ws.Range("N130").Value = -40026.666

# Sheet LTW, row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Sheet LTW, row 92
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 36194
$ws.Range("J92").Value = 36194
$ws.Range("L92").Value = 36194
$ws.Range("N92").Value = -41186

# Sheet LTW, row 128
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 41997.25
$ws.Range("J128").Value = 41997.25
$ws.Range("L128").Value = 41997.25
$ws.Range("N128").Value = -51957.25

# Sheet WVR, row 93
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 28393.4
$ws.Range("J93").Value = 28393.4
$ws.Range("L93").Value = 28393.4
$ws.Range("N93").Value = -33385.4

# Sheet WVR, row 127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 40689.285
$ws.Range("J127").Value = 40689.285
$ws.Range("L127").Value = 40689.285
$ws.Range("N127").Value = -50609.285
